$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply fills by copying format from reference cells with existing styles ---
# (column order within each row, matches visual layout; string VALUES are
#  assigned afterwards in the exact sequence needed to reproduce shared-string order)
$ws.Range("L6").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("L6").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("F17").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("G17").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("H17").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("I17").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("J17").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("K17").PasteSpecial(-4122)

$ws.Range("A2").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("L6").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("F18").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("G18").PasteSpecial(-4122)
$ws.Range("L6").Copy()
$ws.Range("I18").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("J18").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("K18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Assign cell values in the exact order the original edit introduced them ---
# --- (preserves sharedStrings.xml insertion order) ---
$ws.Range("A18").Value = "SD-WIN"
$ws.Range("B18").Value = "PPO use step distance reward + multiply critic lr vs. Random"
$ws.Range("O18").Value = "map*_use_step_dist_train_by_win"
$ws.Range("C17").Value = "绕圈圈，上下碰，才赢了60%"
$ws.Range("C18").Value = "靠左墙慢慢行驶"
$ws.Range("J17").Value = "会撞几次墙"
$ws.Range("J18").Value = "一开始纠缠，路线不稳"
$ws.Range("K18").Value = "贴墙，但是稳"
$ws.Range("E18").Value = "踟蹰or卡墙，只赢了3次"
$ws.Range("E17").Value = "赢了13次，路线很绕"
$ws.Range("G17").Value = "折返"
$ws.Range("H17").Value = "过不了小房间"
$ws.Range("K17").Value = "很流畅，但是到最后没能量"
$ws.Range("F17").Value = "折返并卡墙"
$ws.Range("I18").Value = "偶尔可以"
$ws.Range("F18").Value = "卡墙，很少能过"
$ws.Range("D17").Value = "还行"
$ws.Range("I17").Value = "not a chance"
$ws.Range("D18").Value = "卡墙"
$ws.Range("G18").Value = "卡墙"

# --- Sheet view updates ---
$excel.ActiveWindow.Zoom = 59
$ws.Range("G42").Select()

